$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 173.5887273333334
$ws.Range("H2").Value = 520.7661820000001
$ws.Range("I2").Value = 0.2624583749605043
$ws.Range("J2").Value = 0.2624583749605043
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.400004
$ws.Range("N2").Value = 460.200012
$ws.Range("O2").Value = 0.5184038265041354
$ws.Range("P2").Value = 0.5184038265041354
$ws.Range("Q2").Value = 26628.51146728825
$ws.Range("R2").Value = 239656.6032055942
$ws.Range("S2").Value = 0.1360594258775826
$ws.Range("T2").Value = 0.1360594258775826

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 173.5887273333334
$ws.Range("H3").Value = 520.7661820000001
$ws.Range("I3").Value = 0.2624583749605043
$ws.Range("J3").Value = 0.2624583749605043
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 42.544782
$ws.Range("N3").Value = 127.634346
$ws.Range("O3").Value = 0.1437769048988047
$ws.Range("P3").Value = 0.1437769048988047
$ws.Range("Q3").Value = 7385.294562054109
$ws.Range("R3").Value = 66467.65105848698
$ws.Range("S3").Value = 0.03773545281659127
$ws.Range("T3").Value = 0.03773545281659125

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 173.5887273333334
$ws.Range("H4").Value = 520.7661820000001
$ws.Range("I4").Value = 0.2624583749605043
$ws.Range("J4").Value = 0.2624583749605043
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 67.227361
$ws.Range("N4").Value = 201.682083
$ws.Range("O4").Value = 0.227189832329958
$ws.Range("P4").Value = 0.227189832329958
$ws.Range("Q4").Value = 11669.91203796857
$ws.Range("R4").Value = 105029.2083417171
$ws.Range("S4").Value = 0.05962787420087022
$ws.Range("T4").Value = 0.05962787420087021

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 173.5887273333334
$ws.Range("H5").Value = 520.7661820000001
$ws.Range("I5").Value = 0.2624583749605043
$ws.Range("J5").Value = 0.2624583749605043
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 32.73617033333333
$ws.Range("N5").Value = 98.208511
$ws.Range("O5").Value = 0.1106294362671018
$ws.Range("P5").Value = 0.1106294362671018
$ws.Range("Q5").Value = 5682.630145930557
$ws.Range("R5").Value = 51143.67131337501
$ws.Range("S5").Value = 0.02903562206546022
$ws.Range("T5").Value = 0.02903562206546021

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 141.053299
$ws.Range("H6").Value = 423.159897
$ws.Range("I6").Value = 0.2132662656560029
$ws.Range("J6").Value = 0.2132662656560029
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.400004
$ws.Range("N6").Value = 460.200012
$ws.Range("O6").Value = 0.5184038265041354
$ws.Range("P6").Value = 0.5184038265041354
$ws.Range("Q6").Value = 21637.5766308132
$ws.Range("R6").Value = 194738.1896773188
$ws.Range("S6").Value = 0.1105580481803194
$ws.Range("T6").Value = 0.1105580481803194

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 141.053299
$ws.Range("H7").Value = 423.159897
$ws.Range("I7").Value = 0.2132662656560029
$ws.Range("J7").Value = 0.2132662656560029
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 42.544782
$ws.Range("N7").Value = 127.634346
$ws.Range("O7").Value = 0.1437769048988047
$ws.Range("P7").Value = 0.1437769048988047
$ws.Range("Q7").Value = 6001.081856335818
$ws.Range("R7").Value = 54009.73670702236
$ws.Range("S7").Value = 0.03066276359534636
$ws.Range("T7").Value = 0.03066276359534635

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 141.053299
$ws.Range("H8").Value = 423.159897
$ws.Range("I8").Value = 0.2132662656560029
$ws.Range("J8").Value = 0.2132662656560029
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 67.227361
$ws.Range("N8").Value = 201.682083
$ws.Range("O8").Value = 0.227189832329958
$ws.Range("P8").Value = 0.227189832329958
$ws.Range("Q8").Value = 9482.641052113941
$ws.Range("R8").Value = 85343.76946902546
$ws.Range("S8").Value = 0.04845192713602358
$ws.Range("T8").Value = 0.04845192713602357

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 141.053299
$ws.Range("H9").Value = 423.159897
$ws.Range("I9").Value = 0.2132662656560029
$ws.Range("J9").Value = 0.2132662656560029
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 32.73617033333333
$ws.Range("N9").Value = 98.208511
$ws.Range("O9").Value = 0.1106294362671018
$ws.Range("P9").Value = 0.1106294362671018
$ws.Range("Q9").Value = 4617.544822142597
$ws.Range("R9").Value = 41557.90339928336
$ws.Range("S9").Value = 0.02359352674431358
$ws.Range("T9").Value = 0.02359352674431357

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 322.1880443333333
$ws.Range("H10").Value = 966.564133
$ws.Range("I10").Value = 0.4871338815973437
$ws.Range("J10").Value = 0.4871338815973436
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.400004
$ws.Range("N10").Value = 460.200012
$ws.Range("O10").Value = 0.5184038265041354
$ws.Range("P10").Value = 0.5184038265041354
$ws.Range("Q10").Value = 49423.64728948551
$ws.Range("R10").Value = 444812.8256053696
$ws.Range("S10").Value = 0.2525320682398754
$ws.Range("T10").Value = 0.2525320682398753

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 322.1880443333333
$ws.Range("H11").Value = 966.564133
$ws.Range("I11").Value = 0.4871338815973437
$ws.Range("J11").Value = 0.4871338815973436
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 42.544782
$ws.Range("N11").Value = 127.634346
$ws.Range("O11").Value = 0.1437769048988047
$ws.Range("P11").Value = 0.1437769048988047
$ws.Range("Q11").Value = 13707.420109168
$ws.Range("R11").Value = 123366.780982512
$ws.Range("S11").Value = 0.07003860176740688
$ws.Range("T11").Value = 0.07003860176740687

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 322.1880443333333
$ws.Range("H12").Value = 966.564133
$ws.Range("I12").Value = 0.4871338815973437
$ws.Range("J12").Value = 0.4871338815973436
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 67.227361
$ws.Range("N12").Value = 201.682083
$ws.Range("O12").Value = 0.227189832329958
$ws.Range("P12").Value = 0.227189832329958
$ws.Range("Q12").Value = 21659.851966281
$ws.Range("R12").Value = 194938.667696529
$ws.Range("S12").Value = 0.1106718648823421
$ws.Range("T12").Value = 0.1106718648823421

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 322.1880443333333
$ws.Range("H13").Value = 966.564133
$ws.Range("I13").Value = 0.4871338815973437
$ws.Range("J13").Value = 0.4871338815973436
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 32.73617033333333
$ws.Range("N13").Value = 98.208511
$ws.Range("O13").Value = 0.1106294362671018
$ws.Range("P13").Value = 0.1106294362671018
$ws.Range("Q13").Value = 10547.20269865955
$ws.Range("R13").Value = 94924.82428793596
$ws.Range("S13").Value = 0.05389134670771924
$ws.Range("T13").Value = 0.05389134670771923

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 24.56519766666666
$ws.Range("H14").Value = 73.69559299999999
$ws.Range("I14").Value = 0.03714147778614916
$ws.Range("J14").Value = 0.03714147778614916
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.400004
$ws.Range("N14").Value = 460.200012
$ws.Range("O14").Value = 0.5184038265041354
$ws.Range("P14").Value = 0.5184038265041354
$ws.Range("Q14").Value = 3768.301420327457
$ws.Range("R14").Value = 33914.71278294711
$ws.Range("S14").Value = 0.01925428420635807
$ws.Range("T14").Value = 0.01925428420635807

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 24.56519766666666
$ws.Range("H15").Value = 73.69559299999999
$ws.Range("I15").Value = 0.03714147778614916
$ws.Range("J15").Value = 0.03714147778614916
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 42.544782
$ws.Range("N15").Value = 127.634346
$ws.Range("O15").Value = 0.1437769048988047
$ws.Range("P15").Value = 0.1437769048988047
$ws.Range("Q15").Value = 1045.120979515242
$ws.Range("R15").Value = 9406.088815637177
$ws.Range("S15").Value = 0.005340086719460236
$ws.Range("T15").Value = 0.005340086719460235

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 24.56519766666666
$ws.Range("H16").Value = 73.69559299999999
$ws.Range("I16").Value = 0.03714147778614916
$ws.Range("J16").Value = 0.03714147778614916
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 67.227361
$ws.Range("N16").Value = 201.682083
$ws.Range("O16").Value = 0.227189832329958
$ws.Range("P16").Value = 0.227189832329958
$ws.Range("Q16").Value = 1651.453411573357
$ws.Range("R16").Value = 14863.08070416022
$ws.Range("S16").Value = 0.008438166110722086
$ws.Range("T16").Value = 0.008438166110722085

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 24.56519766666666
$ws.Range("H17").Value = 73.69559299999999
$ws.Range("I17").Value = 0.03714147778614916
$ws.Range("J17").Value = 0.03714147778614916
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 32.73617033333333
$ws.Range("N17").Value = 98.208511
$ws.Range("O17").Value = 0.1106294362671018
$ws.Range("P17").Value = 0.1106294362671018
$ws.Range("Q17").Value = 804.1704950880024
$ws.Range("R17").Value = 7237.534455792022
$ws.Range("S17").Value = 0.004108940749608766
$ws.Range("T17").Value = 0.004108940749608765
